$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header inherits D1's current style (header formatting) along with its
# current text "EDF-L GDPA", which is what the new E1 header should read.
$ws.Range("D1").Copy($ws.Range("E1"))

# D1 header is renamed to "EDF-L HOPA" now that its old text lives in E1.
$ws.Range("D1").Value = "EDF-L HOPA"

# Updated validation-time data for rows 2-21 (columns B-E)
$ws.Range("B2").Value = 0.006705209990032017
$ws.Range("C2").Value = 0.003213256003800779
$ws.Range("D2").Value = 0.01901740193367004
$ws.Range("E2").Value = 0.2974851039052009
$ws.Range("B3").Value = 0.007528988011181355
$ws.Range("C3").Value = 0.004287699987180531
$ws.Range("D3").Value = 0.02039281999692321
$ws.Range("E3").Value = 0.4313589386455715
$ws.Range("B4").Value = 0.01559059995692223
$ws.Range("C4").Value = 0.003921646007802338
$ws.Range("D4").Value = 0.04478444584645331
$ws.Range("E4").Value = 1.240971662830561
$ws.Range("B5").Value = 0.01255253586918116
$ws.Range("C5").Value = 0.004062746013514698
$ws.Range("D5").Value = 0.03804680638946593
$ws.Range("E5").Value = 1.8349620536156
$ws.Range("B6").Value = 0.07532954635564237
$ws.Range("C6").Value = 0.004089272010605782
$ws.Range("D6").Value = 0.2435152406152338
$ws.Range("E6").Value = 3.563859746567905
$ws.Range("B7").Value = 0.01044674599543214
$ws.Range("C7").Value = 0.006859354092739522
$ws.Range("D7").Value = 0.0395663139782846
$ws.Range("E7").Value = 5.903495343886316
$ws.Range("B8").Value = 0.01488286394858733
$ws.Range("C8").Value = 0.006356798009946942
$ws.Range("D8").Value = 0.06007439993321895
$ws.Range("E8").Value = 8.917017914950847
$ws.Range("B9").Value = 0.0160774539783597
$ws.Range("C9").Value = 0.01183786010020412
$ws.Range("D9").Value = 0.07376340391114354
$ws.Range("E9").Value = 8.323454680964351
$ws.Range("B10").Value = 0.04712870353832841
$ws.Range("C10").Value = 0.007191953973378986
$ws.Range("D10").Value = 0.2668349819630385
$ws.Range("E10").Value = 44.01154752969742
$ws.Range("B11").Value = 0.01758691594237462
$ws.Range("C11").Value = 0.007702202001819387
$ws.Range("D11").Value = 1.253325868938118
$ws.Range("E11").Value = 17.12494585752487
$ws.Range("B12").Value = 0.01499161798041314
$ws.Range("C12").Value = 0.006982235978357494
$ws.Range("D12").Value = 4.164473375100642
$ws.Range("E12").Value = 21.03863808274269
$ws.Range("B13").Value = 0.01602638183394447
$ws.Range("C13").Value = 0.009117024013539777
$ws.Range("D13").Value = 3.457636981569231
$ws.Range("E13").Value = 30.58552459359169
$ws.Range("B14").Value = 0.02295655993744731
$ws.Range("C14").Value = 0.007218307995935902
$ws.Range("D14").Value = 11.347207733877
$ws.Range("E14").Value = 33.06925141692162
$ws.Range("B15").Value = 0.01387708600843325
$ws.Range("C15").Value = 0.008675455967895686
$ws.Range("D15").Value = 15.40075955934823
$ws.Range("E15").Value = 43.04863605260849
$ws.Range("B16").Value = 0.01470177997369319
$ws.Range("C16").Value = 0.007691459988709539
$ws.Range("D16").Value = 26.54123341169208
$ws.Range("E16").Value = 49.77835344791412
$ws.Range("B17").Value = 0.0136604980006814
$ws.Range("C17").Value = 0.00799562799045816
$ws.Range("D17").Value = 45.0160842666775
$ws.Range("E17").Value = 62.70932134628296
$ws.Range("B18").Value = 0.01282558211125433
$ws.Range("C18").Value = 0.008138219979591667
$ws.Range("D18").Value = 72.68288652613759
$ws.Range("E18").Value = 85.93666368961334
$ws.Range("B19").Value = 0.01504788205493241
$ws.Range("C19").Value = 0.009719840115867556
$ws.Range("D19").Value = 104.2913051120937
$ws.Range("E19").Value = 135.0216256427765
$ws.Range("B20").Value = 0.01705668007954955
$ws.Range("C20").Value = 0.01102757192915305
$ws.Range("D20").Value = 156.2174384539574
$ws.Range("E20").Value = 205.0270331001282
$ws.Range("B21").Value = 0.02221864799270406
$ws.Range("C21").Value = 0.01178647800348699
$ws.Range("D21").Value = 302.7233975425362
$ws.Range("E21").Value = 418.5600751304626
